# Y4_B2526_Pediatrics_schedule.xlsx update
#
# - Row 2 (the only session kept) is changed to a new group/date/time/duration:
#     B2: PED-A1-1 -> PED-B1-5
#     E2: 10/09/2025 -> 11/09/2025
#     F2: 10:00:00   -> 08:00:00
#     G2: 180         -> 240
# - Rows 3-6 (the other four pediatrics sessions) are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Group code for the remaining session.
$ws.Range("B2").Value = "PED-B1-5"

# The Date column stores plain text like "11/09/2025" (not a real Excel
# date serial) in this workbook, so force the cell to Text before writing
# the new value - otherwise Excel's COM layer auto-parses the
# dd/mm/yyyy-looking string into a date number.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "11/09/2025"
$ws.Range("E2").NumberFormat = "dd/mm/yyyy"

# Start time (stored as text too, e.g. "10:00:00").
$ws.Range("F2").Value = "08:00:00"

# Duration in minutes (numeric).
$ws.Range("G2").Value = 240

# Drop the other four sessions (rows 3-6) entirely, shifting rows up and
# shrinking the used range down to A1:G2.
$ws.Range("A3:G6").EntireRow.Delete()
